$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44763
$ws.Cells.Item(2, 9).Value = "Primera"
$ws.Cells.Item(2, 10).Value = 80
$ws.Cells.Item(2, 11).Value = 17000
$ws.Cells.Item(2, 12).Value = 18000
$ws.Cells.Item(2, 13).Value = 17500
$ws.Cells.Item(2, 16).Value = 972
$ws.Cells.Item(3, 4).Value = 44651
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 60
$ws.Cells.Item(3, 11).Value = 15000
$ws.Cells.Item(3, 12).Value = 16000
$ws.Cells.Item(3, 13).Value = 15500
$ws.Cells.Item(3, 16).Value = 861
$ws.Cells.Item(4, 4).Value = 44813
$ws.Cells.Item(4, 9).Value = "Primera"
$ws.Cells.Item(4, 10).Value = 100
$ws.Cells.Item(4, 11).Value = 14000
$ws.Cells.Item(4, 12).Value = 15000
$ws.Cells.Item(4, 13).Value = 14500
$ws.Cells.Item(4, 16).Value = 806
$ws.Cells.Item(5, 4).Value = 44637
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 100
$ws.Cells.Item(5, 11).Value = 15000
$ws.Cells.Item(5, 12).Value = 16000
$ws.Cells.Item(5, 13).Value = 15500
$ws.Cells.Item(5, 16).Value = 861
$ws.Cells.Item(6, 4).Value = 44649
$ws.Cells.Item(6, 9).Value = "Primera"
$ws.Cells.Item(6, 10).Value = 60
$ws.Cells.Item(6, 11).Value = 15000
$ws.Cells.Item(6, 12).Value = 16000
$ws.Cells.Item(6, 13).Value = 15500
$ws.Cells.Item(6, 16).Value = 861
$ws.Cells.Item(7, 4).Value = 44790
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 60
$ws.Cells.Item(7, 11).Value = 17000
$ws.Cells.Item(7, 12).Value = 18000
$ws.Cells.Item(7, 13).Value = 17500
$ws.Cells.Item(7, 16).Value = 972
$ws.Cells.Item(8, 4).Value = 44754
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 80
$ws.Cells.Item(8, 11).Value = 16000
$ws.Cells.Item(8, 12).Value = 17000
$ws.Cells.Item(8, 13).Value = 16500
$ws.Cells.Item(8, 16).Value = 917
$ws.Cells.Item(9, 4).Value = 44659
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 80
$ws.Cells.Item(9, 11).Value = 15000
$ws.Cells.Item(9, 12).Value = 16000
$ws.Cells.Item(9, 13).Value = 15500
$ws.Cells.Item(9, 16).Value = 861
$ws.Cells.Item(10, 4).Value = 44771
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 60
$ws.Cells.Item(10, 11).Value = 17000
$ws.Cells.Item(10, 12).Value = 18000
$ws.Cells.Item(10, 13).Value = 17500
$ws.Cells.Item(10, 16).Value = 972
$ws.Cells.Item(11, 4).Value = 44804
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 100
$ws.Cells.Item(11, 11).Value = 15000
$ws.Cells.Item(11, 12).Value = 16000
$ws.Cells.Item(11, 13).Value = 15500
$ws.Cells.Item(11, 16).Value = 861
$ws.Cells.Item(12, 4).Value = 44819
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 60
$ws.Cells.Item(12, 11).Value = 15000
$ws.Cells.Item(12, 12).Value = 15000
$ws.Cells.Item(12, 13).Value = 15000
$ws.Cells.Item(12, 16).Value = 833
$ws.Cells.Item(13, 4).Value = 44832
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 60
$ws.Cells.Item(13, 11).Value = 17000
$ws.Cells.Item(13, 12).Value = 17000
$ws.Cells.Item(13, 13).Value = 17000
$ws.Cells.Item(13, 16).Value = 944
$ws.Cells.Item(14, 4).Value = 44384
$ws.Cells.Item(14, 9).Value = "Primera"
$ws.Cells.Item(14, 10).Value = 120
$ws.Cells.Item(14, 11).Value = 17000
$ws.Cells.Item(14, 12).Value = 18000
$ws.Cells.Item(14, 13).Value = 17500
$ws.Cells.Item(14, 16).Value = 972
$ws.Cells.Item(15, 4).Value = 44384
$ws.Cells.Item(15, 9).Value = "Segunda"
$ws.Cells.Item(15, 10).Value = 60
$ws.Cells.Item(15, 11).Value = 15000
$ws.Cells.Item(15, 12).Value = 15000
$ws.Cells.Item(15, 13).Value = 15000
$ws.Cells.Item(15, 16).Value = 833
$ws.Cells.Item(16, 4).Value = 44664
$ws.Cells.Item(16, 9).Value = "Primera"
$ws.Cells.Item(16, 10).Value = 160
$ws.Cells.Item(16, 11).Value = 15000
$ws.Cells.Item(16, 12).Value = 16000
$ws.Cells.Item(16, 13).Value = 15500
$ws.Cells.Item(16, 16).Value = 861
$ws.Cells.Item(17, 4).Value = 44645
$ws.Cells.Item(17, 9).Value = "Primera"
$ws.Cells.Item(17, 10).Value = 60
$ws.Cells.Item(17, 11).Value = 15000
$ws.Cells.Item(17, 12).Value = 16000
$ws.Cells.Item(17, 13).Value = 15500
$ws.Cells.Item(17, 16).Value = 861
$ws.Cells.Item(18, 4).Value = 44761
$ws.Cells.Item(18, 9).Value = "Primera"
$ws.Cells.Item(18, 10).Value = 100
$ws.Cells.Item(18, 11).Value = 17000
$ws.Cells.Item(18, 12).Value = 18000
$ws.Cells.Item(18, 13).Value = 17500
$ws.Cells.Item(18, 16).Value = 972
$ws.Cells.Item(19, 4).Value = 44785
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 80
$ws.Cells.Item(19, 11).Value = 17000
$ws.Cells.Item(19, 12).Value = 18000
$ws.Cells.Item(19, 13).Value = 17500
$ws.Cells.Item(19, 16).Value = 972
$ws.Cells.Item(20, 4).Value = 44799
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 60
$ws.Cells.Item(20, 11).Value = 15000
$ws.Cells.Item(20, 12).Value = 16000
$ws.Cells.Item(20, 13).Value = 15500
$ws.Cells.Item(20, 16).Value = 861
$ws.Cells.Item(21, 4).Value = 44656
$ws.Cells.Item(21, 9).Value = "Primera"
$ws.Cells.Item(21, 10).Value = 100
$ws.Cells.Item(21, 11).Value = 15000
$ws.Cells.Item(21, 12).Value = 16000
$ws.Cells.Item(21, 13).Value = 15500
$ws.Cells.Item(21, 16).Value = 861
$ws.Cells.Item(22, 4).Value = 44818
$ws.Cells.Item(22, 9).Value = "Primera"
$ws.Cells.Item(22, 10).Value = 60
$ws.Cells.Item(22, 11).Value = 15000
$ws.Cells.Item(22, 12).Value = 15000
$ws.Cells.Item(22, 13).Value = 15000
$ws.Cells.Item(22, 16).Value = 833
$ws.Cells.Item(23, 4).Value = 44811
$ws.Cells.Item(23, 9).Value = "Primera"
$ws.Cells.Item(23, 10).Value = 60
$ws.Cells.Item(23, 11).Value = 14000
$ws.Cells.Item(23, 12).Value = 15000
$ws.Cells.Item(23, 13).Value = 14500
$ws.Cells.Item(23, 16).Value = 806
$ws.Cells.Item(24, 4).Value = 44630
$ws.Cells.Item(24, 9).Value = "Primera"
$ws.Cells.Item(24, 10).Value = 60
$ws.Cells.Item(24, 11).Value = 15000
$ws.Cells.Item(24, 12).Value = 16000
$ws.Cells.Item(24, 13).Value = 15500
$ws.Cells.Item(24, 16).Value = 861
$ws.Cells.Item(25, 4).Value = 44658
$ws.Cells.Item(25, 9).Value = "Primera"
$ws.Cells.Item(25, 10).Value = 80
$ws.Cells.Item(25, 11).Value = 15000
$ws.Cells.Item(25, 12).Value = 16000
$ws.Cells.Item(25, 13).Value = 15500
$ws.Cells.Item(25, 16).Value = 861
$ws.Cells.Item(26, 4).Value = 44797
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 80
$ws.Cells.Item(26, 11).Value = 16000
$ws.Cells.Item(26, 12).Value = 17000
$ws.Cells.Item(26, 13).Value = 16500
$ws.Cells.Item(26, 16).Value = 917
$ws.Cells.Item(27, 4).Value = 44775
$ws.Cells.Item(27, 9).Value = "Primera"
$ws.Cells.Item(27, 10).Value = 100
$ws.Cells.Item(27, 11).Value = 17000
$ws.Cells.Item(27, 12).Value = 18000
$ws.Cells.Item(27, 13).Value = 17500
$ws.Cells.Item(27, 16).Value = 972
$ws.Cells.Item(28, 4).Value = 44809
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 60
$ws.Cells.Item(28, 11).Value = 14000
$ws.Cells.Item(28, 12).Value = 15000
$ws.Cells.Item(28, 13).Value = 14500
$ws.Cells.Item(28, 16).Value = 806
$ws.Cells.Item(29, 4).Value = 44635
$ws.Cells.Item(29, 9).Value = "Primera"
$ws.Cells.Item(29, 10).Value = 100
$ws.Cells.Item(29, 11).Value = 15000
$ws.Cells.Item(29, 12).Value = 16000
$ws.Cells.Item(29, 13).Value = 15500
$ws.Cells.Item(29, 16).Value = 861
$ws.Cells.Item(30, 4).Value = 44642
$ws.Cells.Item(30, 9).Value = "Primera"
$ws.Cells.Item(30, 10).Value = 100
$ws.Cells.Item(30, 11).Value = 15000
$ws.Cells.Item(30, 12).Value = 16000
$ws.Cells.Item(30, 13).Value = 15500
$ws.Cells.Item(30, 16).Value = 861
$ws.Cells.Item(31, 4).Value = 44782
$ws.Cells.Item(31, 9).Value = "Primera"
$ws.Cells.Item(31, 10).Value = 120
$ws.Cells.Item(31, 11).Value = 17000
$ws.Cells.Item(31, 12).Value = 18000
$ws.Cells.Item(31, 13).Value = 17500
$ws.Cells.Item(31, 16).Value = 972
$ws.Cells.Item(32, 4).Value = 44839
$ws.Cells.Item(32, 9).Value = "Primera"
$ws.Cells.Item(32, 10).Value = 100
$ws.Cells.Item(32, 11).Value = 17000
$ws.Cells.Item(32, 12).Value = 18000
$ws.Cells.Item(32, 13).Value = 17500
$ws.Cells.Item(32, 16).Value = 972
$ws.Cells.Item(33, 4).Value = 44769
$ws.Cells.Item(33, 9).Value = "Primera"
$ws.Cells.Item(33, 10).Value = 60
$ws.Cells.Item(33, 11).Value = 17000
$ws.Cells.Item(33, 12).Value = 18000
$ws.Cells.Item(33, 13).Value = 17500
$ws.Cells.Item(33, 16).Value = 972
$ws.Cells.Item(34, 4).Value = 44830
$ws.Cells.Item(34, 9).Value = "Primera"
$ws.Cells.Item(34, 10).Value = 60
$ws.Cells.Item(34, 11).Value = 17000
$ws.Cells.Item(34, 12).Value = 17000
$ws.Cells.Item(34, 13).Value = 17000
$ws.Cells.Item(34, 16).Value = 944
$ws.Cells.Item(35, 4).Value = 44628
$ws.Cells.Item(35, 9).Value = "Primera"
$ws.Cells.Item(35, 10).Value = 60
$ws.Cells.Item(35, 11).Value = 15000
$ws.Cells.Item(35, 12).Value = 16000
$ws.Cells.Item(35, 13).Value = 15500
$ws.Cells.Item(35, 16).Value = 861
$ws.Cells.Item(36, 4).Value = 44847
$ws.Cells.Item(36, 9).Value = "Primera"
$ws.Cells.Item(36, 10).Value = 120
$ws.Cells.Item(36, 11).Value = 17000
$ws.Cells.Item(36, 12).Value = 17000
$ws.Cells.Item(36, 13).Value = 17000
$ws.Cells.Item(36, 16).Value = 944
$ws.Cells.Item(37, 4).Value = 44791
$ws.Cells.Item(37, 9).Value = "Primera"
$ws.Cells.Item(37, 10).Value = 80
$ws.Cells.Item(37, 11).Value = 17000
$ws.Cells.Item(37, 12).Value = 18000
$ws.Cells.Item(37, 13).Value = 17500
$ws.Cells.Item(37, 16).Value = 972
